$wb = $excel.ActiveWorkbook

$paramsSheet = $wb.Worksheets.Item("params")
$paramsSheet.Range("A6").Value = "differential_cost"
$paramsSheet.Range("B6").Value = 0
$paramsSheet.Range("C6").Value = "€/m²BGF"

$hullSheet = $wb.Worksheets.Item("thermal_hull")
$hullSheet.Range("G5:G9").Select() | Out-Null

$paramsSheet.Activate()
$paramsSheet.Range("D20").Select() | Out-Null
